$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "risk" column (I7:I16) values to be fully qualified
# Patient.PatientRisk.<LEVEL> enum references instead of bare level names,
# so the table can be called as a template from the API / rule engine.
$ws.Range("I9").Value  = "Patient.PatientRisk.MEDIUM"
$ws.Range("I10").Value = "Patient.PatientRisk.MEDIUM"
$ws.Range("I11").Value = "Patient.PatientRisk.HIGH"
$ws.Range("I12").Value = "Patient.PatientRisk.HIGH"
$ws.Range("I13").Value = "Patient.PatientRisk.HIGH"
$ws.Range("I14").Value = "Patient.PatientRisk.HIGH"
$ws.Range("I15").Value = "Patient.PatientRisk.HIGH"
$ws.Range("I16").Value = "Patient.PatientRisk.HIGH"
$ws.Range("I7").Value  = "Patient.PatientRisk.LOW"
$ws.Range("I8").Value  = "Patient.PatientRisk.LOW"

# Reflect the updated view/selection state left behind in the saved file
# (scrolled down a bit, cell M14 selected instead of N14).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M14").Select()

$wb.Save()
